# Fruta / hortaliza, semanal
# Insert the latest weekly price record for "Zapallo" (Paine, 1a (cosecha))
# at the top of its data block (row 375), pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 375; this shifts rows 375:417 down to 376:418
$ws.Rows("375:375").Insert()

# Populate the new row with this week's record
$ws.Range("A375").Value = 4
$ws.Range("B375").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C375").Value = "Los Lagos"
$ws.Range("D375").Value = 44946
$ws.Range("E375").Value = 10
$ws.Range("F375").Value = 100112045
$ws.Range("G375").Value = "Zapallo"
$ws.Range("H375").Value = "Paine"
$ws.Range("I375").Value = "1a (cosecha)"
$ws.Range("J375").Value = 1200
$ws.Range("K375").Value = 600
$ws.Range("L375").Value = 600
$ws.Range("M375").Value = 600
$ws.Range("N375").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O375").Value = "Región de O'Higgins"
$ws.Range("P375").Value = 600
$ws.Range("Q375").Value = 1
$ws.Range("R375").Value = "Hortaliza"
